$wb = $excel.ActiveWorkbook

# Update the login sheet cell B9 value (single space -> double space)
$loginSheet = $wb.Worksheets.Item("login")
$loginSheet.Range("B9").Value = "  "

# Select D3 on productInfoData before switching away from it
$productInfoSheet = $wb.Worksheets.Item("productInfoData")
$productInfoSheet.Activate()
$productInfoSheet.Range("D3").Select()

# Activate login sheet and select A11
$loginSheet.Activate()
$loginSheet.Range("A11").Select()
